$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New species rows 15-27: column A is always the repeated "spc" label ---
for ($r = 15; $r -le 27; $r++) {
    $ws.Cells.Item($r, 1).Value = "spc"
}

# --- Column B / C (and the occasional G note) in the exact order the
#     original author typed them (reconstructed from shared-string order) ---
$ws.Range("B15").Value = "CO"
$ws.Range("C15").Value = "[C]=O"

$ws.Range("B17").Value = "He"
$ws.Range("C17").Value = "[He]"

$ws.Range("B18").Value = "CH4"
$ws.Range("C18").Value = "C"

$ws.Range("B19").Value = "C2H2"
$ws.Range("B20").Value = "C2H4"
$ws.Range("B22").Value = "CH2O"
$ws.Range("B21").Value = "C2H6"
$ws.Range("B23").Value = "CH3OH"
$ws.Range("B24").Value = "H2O2"

$ws.Range("C19").Value = "C#C"
$ws.Range("C20").Value = "C=C"
$ws.Range("C21").Value = "CC"
$ws.Range("C22").Value = "C=O"
$ws.Range("C24").Value = "OO"

$ws.Range("B25").Value = "CH3CHO"
$ws.Range("C25").Value = "CC=O"

$ws.Range("B16").Value = "CO2"
$ws.Range("C16").Value = "O=C=O"

$ws.Range("B26").Value = "CH2O2"
$ws.Range("G26").Value = "formic acid"
$ws.Range("G25").Value = "acetaldehyde"
$ws.Range("C26").Value = "C(=O)O"

$ws.Range("B27").Value = "HCOOCH3"
$ws.Range("C27").Value = "O=COC"
$ws.Range("G27").Value = "methyl formate"

# --- G1 header, mirrors column A's "name" header ---
$ws.Range("G1").Value = "name"

# --- D / E / F numeric columns for every new row ---
for ($r = 15; $r -le 27; $r++) {
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = 0
    $ws.Cells.Item($r, 6).Value = 0
}

# --- Formatting ---
# C25 picked up the left-aligned style already used elsewhere in the sheet
$ws.Range("C25").HorizontalAlignment = -4131

# G1 / G25 / G26 / G27 use a faint (white, darker 25%) font colour
$ws.Range("G1").Font.ThemeColor = 2
$ws.Range("G25").Font.ThemeColor = 2
$ws.Range("G26").Font.ThemeColor = 2
$ws.Range("G27").Font.ThemeColor = 2

# Column G width
$ws.Columns.Item(7).ColumnWidth = 10

# --- Selection moves to J21, matching the saved session ---
$ws.Range("J21").Select()
